# Update Financials: Balance Sheet figures for FY ending 2018-07-31 (column D)
# Sheet "PANW" on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PANW")

# Net Receivables
$ws.Range("D43").Value = 467000

# Other Current Assets
$ws.Range("D45").Value = 529400

# Total Current Assets
$ws.Range("D46").Value = 4138500

# Property Plant and Equipment
$ws.Range("D48").Value = 546200

# Other Assets
$ws.Range("D52").Value = 533000

# Total Assets
$ws.Range("D54").Value = 5948900

# Other Current Liabilities
$ws.Range("D59").Value = 1501900

# Total Current Liabilities
$ws.Range("D60").Value = 2101700

# Other Liabilities
$ws.Range("D62").Value = 1317200

# Total Liabilities
$ws.Range("D66").Value = 4788600

# Retained Earnings
$ws.Range("D72").Value = -790700

# Total Stockholder Equity
$ws.Range("D76").Value = 1160300
